$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column R (trainer_email): rows 2-12 changed from prameswari.kristal@nutrifood.co.id
# to fahimhadimaula@gmail.com
for ($row = 2; $row -le 12; $row++) {
    $ws.Range("R$row").Value = "fahimhadimaula@gmail.com"
}

# Column W (trainee_email): row 12 changed from prameswari.kristal@nutrifood.co.id
# to fahimhadimaula@gmail.com
$ws.Range("W12").Value = "fahimhadimaula@gmail.com"

# Column AQ (cc_1): rows 2-12 changed from fahimhadimaula@gmail.com
# to fahimhmaula@gmail.com
for ($row = 2; $row -le 12; $row++) {
    $ws.Range("AQ$row").Value = "fahimhmaula@gmail.com"
}
